$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing rows with corrected figures (USD-adjusted calculations) ---

# Row 92 - 2022-05-16
$ws.Cells.Item(92, 6).Value = 165251644
$ws.Cells.Item(92, 7).Value = 171
$ws.Cells.Item(92, 8).Value = 946

# Row 93 - 2022-05-17
$ws.Cells.Item(93, 6).Value = 243601415.95
$ws.Cells.Item(93, 7).Value = 96
$ws.Cells.Item(93, 8).Value = 1365

# Row 587 - 2024-05-27
$ws.Cells.Item(587, 2).Value = 4130
$ws.Cells.Item(587, 4).Value = 4147.95
$ws.Cells.Item(587, 6).Value = 3934231755.95
$ws.Cells.Item(587, 7).Value = 958958
$ws.Cells.Item(587, 8).Value = 4765

# --- Append new rows 589-593 ---
# Helper approach: force column A cells to remain plain text (not auto-converted
# to a date serial number) by applying a text number format before assignment
# and then resetting the cell style back to Normal so no stray formatting is left.

$ws.Cells.Item(589, 1).NumberFormat = "@"
$ws.Cells.Item(589, 1).Value = "2024-05-29"
$ws.Cells.Item(589, 1).Style = "Normal"
$ws.Cells.Item(589, 2).Value = 4163
$ws.Cells.Item(589, 3).Value = 4150
$ws.Cells.Item(589, 4).Value = 4179.6
$ws.Cells.Item(589, 5).Value = 4020
$ws.Cells.Item(589, 6).Value = 6098469470.3
$ws.Cells.Item(589, 7).Value = 1492620
$ws.Cells.Item(589, 8).Value = 5100

$ws.Cells.Item(590, 1).NumberFormat = "@"
$ws.Cells.Item(590, 1).Value = "2024-05-30"
$ws.Cells.Item(590, 1).Style = "Normal"
$ws.Cells.Item(590, 2).Value = 4319
$ws.Cells.Item(590, 3).Value = 4278.15
$ws.Cells.Item(590, 4).Value = 4398
$ws.Cells.Item(590, 5).Value = 4260
$ws.Cells.Item(590, 6).Value = 10611942957.55
$ws.Cells.Item(590, 7).Value = 2461038
$ws.Cells.Item(590, 8).Value = 7665

$ws.Cells.Item(591, 1).NumberFormat = "@"
$ws.Cells.Item(591, 1).Value = "2024-05-31"
$ws.Cells.Item(591, 1).Style = "Normal"
$ws.Cells.Item(591, 2).Value = 4355.15
$ws.Cells.Item(591, 3).Value = 4314.5
$ws.Cells.Item(591, 4).Value = 4415
$ws.Cells.Item(591, 5).Value = 4225
$ws.Cells.Item(591, 6).Value = 6501896343.9
$ws.Cells.Item(591, 7).Value = 1504369
$ws.Cells.Item(591, 8).Value = 7310

$ws.Cells.Item(592, 1).NumberFormat = "@"
$ws.Cells.Item(592, 1).Value = "2024-06-03"
$ws.Cells.Item(592, 1).Style = "Normal"
$ws.Cells.Item(592, 2).Value = 4377
$ws.Cells.Item(592, 3).Value = 4356.1
$ws.Cells.Item(592, 4).Value = 4428.5
$ws.Cells.Item(592, 5).Value = 4315
$ws.Cells.Item(592, 6).Value = 2479919124.5
$ws.Cells.Item(592, 7).Value = 200
$ws.Cells.Item(592, 8).Value = 3484

$ws.Cells.Item(593, 1).NumberFormat = "@"
$ws.Cells.Item(593, 1).Value = "2024-06-04"
$ws.Cells.Item(593, 1).Style = "Normal"
$ws.Cells.Item(593, 2).Value = 4307.8
$ws.Cells.Item(593, 3).Value = 4368
$ws.Cells.Item(593, 4).Value = 4368
$ws.Cells.Item(593, 5).Value = 4260
$ws.Cells.Item(593, 6).Value = 11573945750.65
$ws.Cells.Item(593, 7).Value = 2699459
$ws.Cells.Item(593, 8).Value = 7565

$wb.Save()
